$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 35 (FLTR / Flutter Entertainment / Travel & leisure) - rows below shift up
$ws.Rows.Item(35).Delete()

# Insert a new row before the row that now holds "HIK" (was row 42, now row 41)
# so that it sits after HLMA (Halma plc) and before HIK (Hikma Pharmaceuticals)
$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = "HL"
$ws.Range("B41").Value = "Hargreaves Lansdown"
$ws.Range("C41").Value = "Financial services"
